# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Macroferia Regional de Talca - Chirimoya"
# at the top of the data block (rows 63-64), pushing all existing records
# (old rows 63-124) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 63.
$ws.Range("A63:A64").EntireRow.Insert()

# New row 63 - "Primera" quality, week of 2022-11-09 (serial 44874)
$ws.Cells.Item(63, 1).Value  = 5
$ws.Cells.Item(63, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(63, 3).Value  = "Maule"
$ws.Cells.Item(63, 4).Value  = 44874
$ws.Cells.Item(63, 5).Value  = 7
$ws.Cells.Item(63, 6).Value  = "Fruta"
$ws.Cells.Item(63, 7).Value  = 100107
$ws.Cells.Item(63, 8).Value  = "Otros"
$ws.Cells.Item(63, 9).Value  = 100107002
$ws.Cells.Item(63, 10).Value = "Chirimoya"
$ws.Cells.Item(63, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 200
$ws.Cells.Item(63, 14).Value = 22000
$ws.Cells.Item(63, 15).Value = 22000
$ws.Cells.Item(63, 16).Value = 22000
$ws.Cells.Item(63, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(63, 19).Value = 2200
$ws.Cells.Item(63, 20).Value = 10

# New row 64 - "Segunda" quality, same week (serial 44874)
$ws.Cells.Item(64, 1).Value  = 5
$ws.Cells.Item(64, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(64, 3).Value  = "Maule"
$ws.Cells.Item(64, 4).Value  = 44874
$ws.Cells.Item(64, 5).Value  = 7
$ws.Cells.Item(64, 6).Value  = "Fruta"
$ws.Cells.Item(64, 7).Value  = 100107
$ws.Cells.Item(64, 8).Value  = "Otros"
$ws.Cells.Item(64, 9).Value  = 100107002
$ws.Cells.Item(64, 10).Value = "Chirimoya"
$ws.Cells.Item(64, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(64, 12).Value = "Segunda"
$ws.Cells.Item(64, 13).Value = 150
$ws.Cells.Item(64, 14).Value = 20000
$ws.Cells.Item(64, 15).Value = 20000
$ws.Cells.Item(64, 16).Value = 20000
$ws.Cells.Item(64, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(64, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 19).Value = 2000
$ws.Cells.Item(64, 20).Value = 10
